$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append new row 75 with data for 2020-08-13
# Force the date-like text to stay as plain text (matches existing column A strings),
# then restore the default (Normal) style so no explicit style is left on the cell.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "2020-08-13"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = 505751
$ws.Range("C75").Value = 553219
$ws.Range("D75").Value = 83075
$ws.Range("E75").Value = 55293
$ws.Range("F75").Value = 26.39
